# Update Betfair back/lay odds values on Sheet1 (rows 2-15)
# to reflect the latest odds scrape for 2025-11-20.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 970
$ws.Range("I2").Value = 970
$ws.Range("J2").Value = 1.03
$ws.Range("O2").Value = 1.01
$ws.Range("F3").Value = 3.05
$ws.Range("G3").Value = 3.5
$ws.Range("H3").Value = 2.22
$ws.Range("I3").Value = 2.46
$ws.Range("J3").Value = 3.55
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 4.3
$ws.Range("Q3").Value = 1.72
$ws.Range("R3").Value = 1.44
$ws.Range("S3").Value = 2.8
$ws.Range("T3").Value = 1.61
$ws.Range("U3").Value = 2.32
$ws.Range("V3").Value = 1.68
$ws.Range("W3").Value = 1.4
$ws.Range("X3").Value = 23
$ws.Range("Y3").Value = 13
$ws.Range("AA3").Value = 980
$ws.Range("AB3").Value = 16
$ws.Range("AC3").Value = 10.5
$ws.Range("AD3").Value = 14
$ws.Range("AE3").Value = 24
$ws.Range("AF3").Value = 26
$ws.Range("AG3").Value = 15
$ws.Range("AH3").Value = 17
$ws.Range("AI3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 48
$ws.Range("AN3").Value = 34
$ws.Range("AO3").Value = 18
$ws.Range("F4").Value = 2.72
$ws.Range("G4").Value = 3.6
$ws.Range("H4").Value = 2.52
$ws.Range("I4").Value = 3.35
$ws.Range("J4").Value = 2.38
$ws.Range("K4").Value = 4.1
$ws.Range("L4").Value = 1.4
$ws.Range("M4").Value = 1.09
$ws.Range("N4").Value = 2.7
$ws.Range("O4").Value = 1.39
$ws.Range("P4").Value = 1.66
$ws.Range("Q4").Value = 2.02
$ws.Range("R4").Value = 1.25
$ws.Range("S4").Value = 3.45
$ws.Range("T4").Value = 1.83
$ws.Range("U4").Value = 1.9
$ws.Range("V4").Value = 1.43
$ws.Range("W4").Value = 1.38
$ws.Range("F5").Value = 1.64
$ws.Range("G5").Value = 1.78
$ws.Range("H5").Value = 7
$ws.Range("I5").Value = 8.800000000000001
$ws.Range("J5").Value = 3.3
$ws.Range("K5").Value = 3.8
$ws.Range("L5").Value = 1.51
$ws.Range("M5").Value = 1.13
$ws.Range("N5").Value = 2.42
$ws.Range("O5").Value = 1.56
$ws.Range("P5").Value = 1.47
$ws.Range("Q5").Value = 2.7
$ws.Range("R5").Value = 1.16
$ws.Range("S5").Value = 5.1
$ws.Range("T5").Value = 2.5
$ws.Range("U5").Value = 1.54
$ws.Range("V5").Value = 1.13
$ws.Range("W5").Value = 2.28
$ws.Range("X5").Value = 10
$ws.Range("Y5").Value = 20
$ws.Range("AB5").Value = 6.4
$ws.Range("AC5").Value = 10.5
$ws.Range("AD5").Value = 42
$ws.Range("AF5").Value = 9.800000000000001
$ws.Range("AG5").Value = 14
$ws.Range("AH5").Value = 44
$ws.Range("AJ5").Value = 22
$ws.Range("AK5").Value = 32
$ws.Range("AN5").Value = 24
$ws.Range("G6").Value = 5.1
$ws.Range("I6").Value = 1.9
$ws.Range("J6").Value = 3.85
$ws.Range("L6").Value = 1.45
$ws.Range("P6").Value = 1.79
$ws.Range("R6").Value = 1.3
$ws.Range("V6").Value = 2.1
$ws.Range("W6").Value = 1.25
$ws.Range("Z6").Value = 11
$ws.Range("AL6").Value = 80
$ws.Range("AN6").Value = 90
$ws.Range("F7").Value = 1.33
$ws.Range("G7").Value = 1.59
$ws.Range("K7").Value = 8
$ws.Range("W7").Value = 2.68
$ws.Range("H8").Value = 4.9
$ws.Range("P8").Value = 1.88
$ws.Range("Q8").Value = 1.92
$ws.Range("F9").Value = 1.54
$ws.Range("I9").Value = 7.2
$ws.Range("J9").Value = 4.7
$ws.Range("K9").Value = 4.9
$ws.Range("L9").Value = 1.32
$ws.Range("M9").Value = 1.04
$ws.Range("N9").Value = 4.8
$ws.Range("O9").Value = 1.22
$ws.Range("P9").Value = 2.28
$ws.Range("Q9").Value = 1.67
$ws.Range("R9").Value = 1.51
$ws.Range("S9").Value = 2.64
$ws.Range("T9").Value = 1.78
$ws.Range("U9").Value = 2.12
$ws.Range("V9").Value = 1.16
$ws.Range("Y9").Value = 25
$ws.Range("Z9").Value = 60
$ws.Range("AB9").Value = 10.5
$ws.Range("AD9").Value = 25
$ws.Range("AI9").Value = 75
$ws.Range("AJ9").Value = 15
$ws.Range("AM9").Value = 110
$ws.Range("F10").Value = 1.66
$ws.Range("G10").Value = 2.14
$ws.Range("H10").Value = 4.2
$ws.Range("K10").Value = 6.2
$ws.Range("N10").Value = 1.94
$ws.Range("W10").Value = 1.89
$ws.Range("I11").Value = 3.75
$ws.Range("K11").Value = 3.15
$ws.Range("L11").Value = 1.57
$ws.Range("N11").Value = 2.64
$ws.Range("S11").Value = 5.3
$ws.Range("T11").Value = 2.06
$ws.Range("AB11").Value = 7.8
$ws.Range("AM11").Value = 200
$ws.Range("F12").Value = 1.99
$ws.Range("G12").Value = 2.14
$ws.Range("H12").Value = 4.2
$ws.Range("I12").Value = 5.2
$ws.Range("P12").Value = 1.56
$ws.Range("Q12").Value = 2.04
$ws.Range("G13").Value = 2.36
$ws.Range("H13").Value = 3.75
$ws.Range("P13").Value = 1.55
$ws.Range("S13").Value = 5.4
$ws.Range("Z13").Value = 27
$ws.Range("AD13").Value = 970
$ws.Range("AF13").Value = 970
$ws.Range("AG13").Value = 970
$ws.Range("AH13").Value = 24
$ws.Range("AJ13").Value = 34
$ws.Range("AK13").Value = 32
$ws.Range("F14").Value = 1.68
$ws.Range("I14").Value = 6.8
$ws.Range("J14").Value = 3.8
$ws.Range("K14").Value = 4.3
$ws.Range("M14").Value = 1.07
$ws.Range("N14").Value = 3.55
$ws.Range("P14").Value = 1.9
$ws.Range("Q14").Value = 1.96
$ws.Range("R14").Value = 1.34
$ws.Range("T14").Value = 1.77
$ws.Range("U14").Value = 1.92
$ws.Range("V14").Value = 1.19
$ws.Range("AC14").Value = 11.5
$ws.Range("AD14").Value = 27
$ws.Range("AE14").Value = 100
$ws.Range("AI14").Value = 100
$ws.Range("G15").Value = 2.16
$ws.Range("I15").Value = 4
$ws.Range("M15").Value = 1.05
$ws.Range("N15").Value = 4.2
$ws.Range("R15").Value = 1.42
$ws.Range("S15").Value = 2.9
$ws.Range("T15").Value = 1.66
$ws.Range("U15").Value = 2.2
$ws.Range("V15").Value = 1.33
$ws.Range("W15").Value = 1.86
$ws.Range("X15").Value = 18.5
$ws.Range("Y15").Value = 17
$ws.Range("AB15").Value = 11.5
$ws.Range("AC15").Value = 9.4
$ws.Range("AD15").Value = 16.5
$ws.Range("AF15").Value = 14.5
$ws.Range("AG15").Value = 11.5
$ws.Range("AH15").Value = 18
$ws.Range("AJ15").Value = 26
$ws.Range("AK15").Value = 22
$ws.Range("AL15").Value = 34
$ws.Range("AN15").Value = 14
$ws.Range("AO15").Value = 44
